$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
# G48 was never used in the original file; let's touch it and clear to see if it appears as stub.
$ws.Range("G48").Value = "temp"
$ws.Range("G48").ClearContents()
Write-Host "done"
